$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells we touch to remain plain text so that
# numeric-looking values (e.g. "0.540", "379.30") keep their exact original
# formatting instead of being auto-coerced into numbers by Excel.
$priceCells = @("D2","D3","D4","D5","D6","D7","D10","D12","D13","D14","D15","D16","D17","D18","D19","D21","D23","D24","D29","D32","D33","D35","D38","D40","D43","D44","D45","D47","D48","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated cryptocurrency prices
$ws.Range("D2").Value = "51.004.46"
$ws.Range("D3").Value = "2.949.68"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "379.30"
$ws.Range("D6").Value = "101.06"
$ws.Range("D7").Value = "0.540"
$ws.Range("D10").Value = "36.17"
$ws.Range("D12").Value = "0.0850"
$ws.Range("D13").Value = "3.406.42"
$ws.Range("D14").Value = "18.29"
$ws.Range("D15").Value = "7.71"
$ws.Range("D16").Value = "11.89"
$ws.Range("D17").Value = "2.939.73"
$ws.Range("D18").Value = "0.997"
$ws.Range("D19").Value = "50.938.72"
$ws.Range("D21").Value = "12.38"
$ws.Range("D23").Value = "69.49"
$ws.Range("D24").Value = "266.84"
$ws.Range("D29").Value = "25.65"
$ws.Range("D32").Value = "10.09"
$ws.Range("D33").Value = "50.48"
$ws.Range("D35").Value = "33.48"
$ws.Range("D38").Value = "3.11"
$ws.Range("D40").Value = "16.58"
$ws.Range("D43").Value = "120.30"
$ws.Range("D44").Value = "21.42"
$ws.Range("D45").Value = "3.46"
$ws.Range("D47").Value = "2.32"
$ws.Range("D48").Value = "2.011.24"
$ws.Range("D50").Value = "0.0314"

# Restore the default (Normal) style on the price cells so only their
# textual content changed, matching the original unstyled cells.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Updated 1h volume/change percentages
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("E16").Value = "  +66.78%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -5.01%  "
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("E25").Value = "  +13.14%  "
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -6.95%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +4.61%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +7.24%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("E51").Value = "  +4.53%  "
